$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Define the print area for the sheet (A1:G14) - becomes the
# _xlnm.Print_Area defined name scoped to this sheet.
$ws.PageSetup.PrintArea = '$A$1:$G$14'

# Page setup: paper size (9 = A4) and portrait orientation, plus the
# page margins (values are in points; 0.7874in == 56.6928pt, 0.315in == 22.68pt).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.LeftMargin = 56.6928
$ws.PageSetup.RightMargin = 56.6928
$ws.PageSetup.TopMargin = 56.6928
$ws.PageSetup.BottomMargin = 56.6928
$ws.PageSetup.HeaderMargin = 22.68
$ws.PageSetup.FooterMargin = 22.68

# Best-fit the data columns (A:G) to their content, matching the
# bestFit column widths baked into the authored workbook.
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws.Columns.Item(5).ColumnWidth = 17.0
$ws.Columns.Item(6).ColumnWidth = 17.0
$ws.Columns.Item(7).ColumnWidth = 11.833333333333334
